# Add a "Phone number" field (column D) to the file-history bulk upload
# template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in D1 ("Phone number"); Excel auto-adds the shared string,
# extends the used-range dimension to A1:D6, inherits the bold header
# style (s=2) from the row, and extends every row's "spans" to 1:4.
$ws.Range("D1").Value = "Phone number"

# Give column D an explicit width, matching the other header columns.
# (ColumnWidth is quantized to whole pixels by Excel, same as the native
# app, so this lands on the nearest achievable width to 21.7265625 chars.)
$ws.Columns.Item(4).ColumnWidth = 20.8333333333333

# Move the active selection to E2 (one column past the new field), as in
# the source workbook.
[void]$ws.Range("E2").Select()
